$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.963.02"
$ws.Range("D3").Value = "3.672.77"
$ws.Range("E3").Value = "  +18.36%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'619.30"
$ws.Range("E5").Value = "  +7.44%  "
$ws.Range("D6").Value = "'182.50"
$ws.Range("E6").Value = "  +2.91%  "
$ws.Range("D7").Value = "3.669.91"
$ws.Range("E7").Value = "  +18.36%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +5.70%  "
$ws.Range("E10").Value = "  +8.10%  "
$ws.Range("E11").Value = "  +5.37%  "
$ws.Range("E12").Value = "  +7.19%  "
$ws.Range("D13").Value = "'40.41"
$ws.Range("E13").Value = "  +11.85%  "
$ws.Range("D15").Value = "4.281.33"
$ws.Range("E15").Value = "  +18.21%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.675.18"
$ws.Range("E16").Value = "  +18.47%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "70.985.51"
$ws.Range("E17").Value = "  +6.10%  "
$ws.Range("E18").Value = "  +2.05%  "
$ws.Range("E19").Value = "  +7.35%  "
$ws.Range("D20").Value = "'520.51"
$ws.Range("E20").Value = "  +8.50%  "
$ws.Range("D21").Value = "'16.93"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").Value = "'9.26"
$ws.Range("E22").Value = "  +18.92%  "
$ws.Range("D23").Value = "'0.743"
$ws.Range("E23").Value = "  +7.68%  "
$ws.Range("E24").Value = "  +12.80%  "
$ws.Range("D25").Value = "'88.57"
$ws.Range("E25").Value = "  +5.97%  "
$ws.Range("D26").Value = "'13.52"
$ws.Range("E26").Value = "  +7.59%  "
$ws.Range("E27").Value = "  +9.48%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "'2.54"
$ws.Range("E29").Value = "  +10.66%  "
$ws.Range("D30").Value = "'8.18"
$ws.Range("E30").Value = "  +3.51%  "
$ws.Range("D31").Value = "'2.92"
$ws.Range("E32").Value = "  +18.19%  "
$ws.Range("D33").Value = "'31.64"
$ws.Range("E33").Value = "  +12.96%  "
$ws.Range("E34").Value = "  +4.46%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'6.14"
$ws.Range("E36").Value = "  +9.71%  "
$ws.Range("E37").Value = "  +9.30%  "
$ws.Range("D38").Value = "'0.347"
$ws.Range("E38").Value = "  +11.23%  "
$ws.Range("D39").Value = "'2.22"
$ws.Range("E39").Value = "  +9.80%  "
$ws.Range("E40").Value = "  +6.89%  "
$ws.Range("D41").Value = "'51.45"
$ws.Range("E41").Value = "  +4.85%  "
$ws.Range("D42").Value = "'45.64"
$ws.Range("E42").Value = "  -5.72%  "
$ws.Range("D43").Value = "'432.96"
$ws.Range("E43").Value = "  +16.17%  "
$ws.Range("E44").Value = "  +6.00%  "
$ws.Range("D45").Value = "3.110.68"
$ws.Range("E45").Value = "  +11.24%  "
$ws.Range("D46").Value = "'2.84"
$ws.Range("E46").Value = "  +4.90%  "
$ws.Range("E47").Value = "  +7.11%  "
$ws.Range("E48").Value = "  +9.80%  "
$ws.Range("D49").Value = "'140.26"
$ws.Range("E49").Value = "  +3.27%  "
$ws.Range("E51").Value = "  +10.12%  "
